# edit.ps1 -- "New crime data collected" weekly refresh for the CompStat sheet.
#
# Source report moves forward one week (Volume 32, Number 16 -> 17; the
# covering week becomes 4/21/2025-4/27/2025) and every crime-category row in
# the Week/28-Day/YTD/2-Year table (rows 14-31) is refreshed with the new
# weekly counts + recomputed % change figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: bump the volume/issue number and the covering-week dates ---
$ws.Range("A8").Value = 'Volume 32   Number  17'
$ws.Range("C9").Value = 'Report Covering the Week  4/21/2025  Through  4/27/2025'

# --- Crime Complaints table (Week to Date / 28 Day / YTD / 2 Year cols) ---
# A couple of cells (C14, F14, C31) were previously blank ("N/A", shown as the
# shared text "-") and now carry real counts, so their number format is
# switched from the text style to the "#,##0" numeric style used by the rest
# of the column.

# Row 14 - Murder
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("F14").Value = 1
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = -33.333333333333
$ws.Range("L14").Value = 33.333333333333
$ws.Range("M14").Value = -42.857142857142
$ws.Range("N14").Value = -87.5

# Row 15 - Rape
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 700
$ws.Range("F15").Value = 22
$ws.Range("G15").Value = 7
$ws.Range("H15").Value = 214.285714285714
$ws.Range("I15").Value = 69
$ws.Range("J15").Value = 45
$ws.Range("K15").Value = 53.333333333333
$ws.Range("L15").Value = 109.090909090909
$ws.Range("M15").Value = 130
$ws.Range("N15").Value = 9.523809523809

# Row 16 - Robbery
$ws.Range("C16").Value = 38
$ws.Range("D16").Value = 34
$ws.Range("E16").Value = 11.764705882352
$ws.Range("F16").Value = 134
$ws.Range("G16").Value = 136
$ws.Range("H16").Value = -1.470588235294
$ws.Range("I16").Value = 471
$ws.Range("J16").Value = 486
$ws.Range("K16").Value = -3.086419753086
$ws.Range("L16").Value = -16.931216931216
$ws.Range("M16").Value = 29.752066115702
$ws.Range("N16").Value = -85.156003781909

# Row 17 - Fel. Assault
$ws.Range("C17").Value = 46
$ws.Range("D17").Value = 51
$ws.Range("E17").Value = -9.803921568627
$ws.Range("F17").Value = 197
$ws.Range("G17").Value = 176
$ws.Range("H17").Value = 11.931818181818
$ws.Range("I17").Value = 685
$ws.Range("J17").Value = 683
$ws.Range("K17").Value = 0.292825768667
$ws.Range("L17").Value = 0.292825768667
$ws.Range("M17").Value = 72.110552763819
$ws.Range("N17").Value = -32.843137254902

# Row 18 - Burglary
$ws.Range("C18").Value = 29
$ws.Range("D18").Value = 54
$ws.Range("E18").Value = -46.296296296296
$ws.Range("F18").Value = 138
$ws.Range("G18").Value = 145
$ws.Range("H18").Value = -4.827586206896
$ws.Range("I18").Value = 639
$ws.Range("J18").Value = 624
$ws.Range("K18").Value = 2.403846153846
$ws.Range("L18").Value = -7.658959537572
$ws.Range("M18").Value = 2.733118971061
$ws.Range("N18").Value = -83.014354066985

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 163
$ws.Range("D19").Value = 194
$ws.Range("E19").Value = -15.979381443299
$ws.Range("F19").Value = 668
$ws.Range("G19").Value = 721
$ws.Range("H19").Value = -7.350901525658
$ws.Range("I19").Value = 3038
$ws.Range("J19").Value = 3189
$ws.Range("K19").Value = -4.735026654123
$ws.Range("L19").Value = -15.681376630585
$ws.Range("M19").Value = -4.973412574288
$ws.Range("N19").Value = -69.328621908127

# Row 20 - G.L.A.
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 12.5
$ws.Range("F20").Value = 26
$ws.Range("H20").Value = -10.344827586206
$ws.Range("I20").Value = 91
$ws.Range("J20").Value = 126
$ws.Range("K20").Value = -27.777777777777
$ws.Range("L20").Value = -45.833333333333
$ws.Range("M20").Value = -15.74074074074
$ws.Range("N20").Value = -95.447723861931

# Row 21 - TOTAL
$ws.Range("C21").Value = 294
$ws.Range("D21").Value = 342
$ws.Range("E21").Value = -14.035087719298
$ws.Range("F21").Value = 1186
$ws.Range("G21").Value = 1214
$ws.Range("H21").Value = -2.306425041186
$ws.Range("I21").Value = 4997
$ws.Range("J21").Value = 5159
$ws.Range("K21").Value = -3.14014343865
$ws.Range("L21").Value = -13.080535745347
$ws.Range("M21").Value = 5.756613756613
$ws.Range("N21").Value = -74.957402024656

# Row 22 - Transit
$ws.Range("C22").Value = 9
$ws.Range("D22").Value = 9
$ws.Range("F22").Value = 42
$ws.Range("G22").Value = 47
$ws.Range("H22").Value = -10.63829787234
$ws.Range("I22").Value = 212
$ws.Range("J22").Value = 204
$ws.Range("K22").Value = 3.92156862745
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 22.543352601156

# Row 23 - Housing
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 13
$ws.Range("E23").Value = -30.76923076923
$ws.Range("F23").Value = 33
$ws.Range("G23").Value = 35
$ws.Range("H23").Value = -5.714285714285
$ws.Range("I23").Value = 140
$ws.Range("J23").Value = 124
$ws.Range("K23").Value = 12.903225806451
$ws.Range("L23").Value = 14.754098360655
$ws.Range("M23").Value = 20.689655172413

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 379
$ws.Range("D24").Value = 446
$ws.Range("E24").Value = -15.022421524663
$ws.Range("F24").Value = 1536
$ws.Range("G24").Value = 1615
$ws.Range("H24").Value = -4.891640866873
$ws.Range("I24").Value = 6339
$ws.Range("J24").Value = 6785
$ws.Range("K24").Value = -6.573323507737
$ws.Range("L24").Value = 5.526885300482
$ws.Range("M24").Value = 25.948738327041

# Row 25 - Retail Theft
$ws.Range("C25").Value = 284
$ws.Range("D25").Value = 391
$ws.Range("E25").Value = -27.365728900255
$ws.Range("F25").Value = 1152
$ws.Range("G25").Value = 1391
$ws.Range("H25").Value = -17.181883537023
$ws.Range("I25").Value = 5135
$ws.Range("J25").Value = 5821
$ws.Range("K25").Value = -11.784916680982
$ws.Range("L25").Value = 1.382033563672

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 90
$ws.Range("D26").Value = 91
$ws.Range("E26").Value = -1.098901098901
$ws.Range("F26").Value = 372
$ws.Range("G26").Value = 383
$ws.Range("H26").Value = -2.872062663185
$ws.Range("I26").Value = 1536
$ws.Range("J26").Value = 1532
$ws.Range("K26").Value = 0.261096605744
$ws.Range("L26").Value = 4.135593220338
$ws.Range("M26").Value = 32.299741602067

# Row 27 - UCR Rape*
$ws.Range("C27").Value = 9
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 350
$ws.Range("F27").Value = 24
$ws.Range("G27").Value = 16
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 81
$ws.Range("J27").Value = 72
$ws.Range("K27").Value = 12.5
$ws.Range("L27").Value = 19.117647058823

# Row 28 - Other Sex Crimes
$ws.Range("C28").Value = 23
$ws.Range("D28").Value = 18
$ws.Range("E28").Value = 27.777777777777
$ws.Range("F28").Value = 78
$ws.Range("G28").Value = 80
$ws.Range("H28").Value = -2.5
$ws.Range("I28").Value = 280
$ws.Range("J28").Value = 278
$ws.Range("K28").Value = 0.719424460431
$ws.Range("L28").Value = -0.709219858156

# Row 29 - Shooting Vic.
$ws.Range("G29").Value = 2
$ws.Range("M29").Value = -75
$ws.Range("N29").Value = -94.230769230769

# Row 30 - Shooting Inc.
$ws.Range("G30").Value = 2
$ws.Range("M30").Value = -62.5
$ws.Range("N30").Value = -93.333333333333

# Row 31 - Hate Crimes
$ws.Range("C31").Value = 2
$ws.Range("C31").NumberFormat = "#,##0"
$ws.Range("E31").Value = -50
$ws.Range("F31").Value = 15
$ws.Range("H31").Value = -16.666666666666
$ws.Range("I31").Value = 47
$ws.Range("J31").Value = 49
$ws.Range("K31").Value = -4.081632653061
$ws.Range("L31").Value = 9.302325581395
